$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 4
$ws.Range("F7").Value = "牟秋宇"
$ws.Range("B7").Value = "完成非功能性测试中性能、可维护性、可访问性、安全性、界面美观等测试"

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "完成测试规格说明书中非功能性测试部分"

$ws.Range("E7").Value = "有"
$ws.Range("F8").Value = "牟秋宇"

$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 28.8

$null = $ws.Range("G7").Select()
